# Updates scraped cryptocurrency prices / 1h volume changes (cryptos list refresh).
# A handful of rows also swap rank position (Bittensor <-> RenderToken, SuiNetwork <-> Filecoin).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "60.130.10"; E = "  +2.46%  " },
    @{ Row = 3; D = "2.549.00"; E = "  +1.48%  " },
    @{ Row = 4; E = "  -0.01%  " },
    @{ Row = 5; D = "540.25"; E = "  +1.47%  " },
    @{ Row = 6; D = "144.18"; E = "  +1.55%  " },
    @{ Row = 7; D = "0.997"; E = "  -0.21%  " },
    @{ Row = 8; E = "  +0.21%  " },
    @{ Row = 9; D = "2.569.31"; E = "  +2.30%  " },
    @{ Row = 10; E = "  +2.00%  " },
    @{ Row = 11; E = "  +1.82%  " },
    @{ Row = 12; D = "5.48"; E = "  +1.06%  " },
    @{ Row = 13; E = "  +3.82%  " },
    @{ Row = 14; D = "2.997.58"; E = "  +1.55%  " },
    @{ Row = 15; D = "24.09"; E = "  +1.57%  " },
    @{ Row = 16; D = "60.073.16"; E = "  +2.40%  " },
    @{ Row = 17; D = "0.0000143"; E = "  +4.44%  " },
    @{ Row = 18; D = "2.560.33"; E = "  +2.43%  " },
    @{ Row = 19; E = "  -1.01%  " },
    @{ Row = 20; E = "  +1.60%  " },
    @{ Row = 21; D = "327.03"; E = "  +1.59%  " },
    @{ Row = 22; D = "0.999"; E = "  +0.03%  " },
    @{ Row = 23; D = "5.96"; E = "  +3.83%  " },
    @{ Row = 24; D = "63.33"; E = "  +4.31%  " },
    @{ Row = 25; E = "  -0.68%  " },
    @{ Row = 26; D = "0.167"; E = "  +4.02%  " },
    @{ Row = 27; E = "  -0.01%  " },
    @{ Row = 28; D = "8.04"; E = "  +4.28%  " },
    @{ Row = 29; E = "  +3.23%  " },
    @{ Row = 30; D = "0.0₃0796"; E = "  +4.42%  " },
    @{ Row = 31; E = "  +2.39%  " },
    @{ Row = 32; E = "  -3.54%  " },
    @{ Row = 33; D = "165.47"; E = "  +5.39%  " },
    @{ Row = 34; E = "  +5.33%  " },
    @{ Row = 35; E = "  +0.12%  " },
    @{ Row = 36; D = "18.74"; E = "  +1.30%  " },
    @{ Row = 37; E = "  +0.93%  " },
    @{ Row = 38; E = "  +2.51%  " },
    @{ Row = 39; D = "37.04"; E = "  +0.90%  " },
    @{ Row = 40; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "5.61"; E = "  -5.25%  " },
    @{ Row = 41; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "301.87"; E = "  -1.40%  " },
    @{ Row = 42; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "3.72"; E = "  +2.16%  " },
    @{ Row = 43; B = "SuiNetwork"; C = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D = "0.834"; E = "  +5.96%  " },
    @{ Row = 44; D = "0.612"; E = "  +2.98%  " },
    @{ Row = 45; E = "  -0.48%  " },
    @{ Row = 46; E = "  +0.94%  " },
    @{ Row = 47; D = "127.21"; E = "  +2.55%  " },
    @{ Row = 48; D = "0.0939"; E = "  +1.70%  " },
    @{ Row = 49; E = "  +1.97%  " },
    @{ Row = 50; E = "  +0.92%  " },
    @{ Row = 51; D = "0.0229"; E = "  +1.29%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($r, 3).Value = $u.C }

    if ($u.ContainsKey("D")) {
        # Price column holds text (e.g. "60.130.10", "0.999"). Force text entry so Excel
        # doesn't silently reinterpret number-looking values as numerics, then drop back
        # to the default style so no formatting residue is left on the cell.
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }

    if ($u.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $u.E }
}
